# Apply the cryptos-list refresh described by the commit diff.
# D-column "Price" cells are free-form text (often contain dotted
# thousands separators, e.g. "27.394.97"), and E-column "Volume(1h)"
# cells are padded percentage strings. Excel's COM layer auto-coerces a
# plain numeric-looking string (e.g. "20.05") into a real number when
# assigned straight to .Value, so for D-column values that parse as a
# bare number we lead with an apostrophe (forces text entry, same as
# typing it in the Excel UI) and then reset the cell Style back to
# "Normal" so no stray quote-prefix / number-format style sticks to the
# cell (matching the source diff, which only changes the text content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.394.97'
$ws.Range("D3").Value = '1.861.93'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = '''315.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").Value = '''1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").Value = '''0.4622'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").Value = '''0.3721'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").Value = '''0.07324'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").Value = '''0.8896'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.10%  '
$ws.Range("D11").Value = '''20.05'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '''0.07853'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").Value = '1.922.83'
$ws.Range("E13").Value = '  +6.26%  '
$ws.Range("D14").Value = '''5.399'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '''6.556'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = '''91.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '''0.000008974'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").Value = '''1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '''14.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("D21").Value = '27.399.45'
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("D22").Value = '''5.134'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '''10.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '2.108.58'
$ws.Range("E24").Value = '  +5.28%  '
$ws.Range("D25").Value = '''1.937'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.26%  '
$ws.Range("D26").Value = '''152.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = '''18.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '''2.054'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").Value = '''5.103'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").Value = '''116.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = '''3.089'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.86%  '
$ws.Range("D33").Value = '''0.7670'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.70%  '
$ws.Range("D34").Value = '''1.177'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").Value = '''4.521'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.93%  '
$ws.Range("D36").Value = '''2.712'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.42%  '
$ws.Range("D37").Value = '''1.081'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("D38").Value = '''0.01959'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("D39").Value = '''2.985'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("D40").Value = '''0.05247'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").Value = '''7.081'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").Value = '''0.5152'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").Value = '''8.426'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.36%  '
$ws.Range("D45").Value = '''0.4812'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '''10.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").Value = '''1.003'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").Value = '''102.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("D49").Value = '''1.648'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.37%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").Value = '''65.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.97%  '
